$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 27.90240433333333
$ws.Range("H2").Value = 83.707213
$ws.Range("I2").Value = 0.2174736967445081
$ws.Range("J2").Value = 0.2174736967445081
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2296153333333334
$ws.Range("N2").Value = 0.6888460000000001
$ws.Range("O2").Value = 0.04381656765025366
$ws.Range("P2").Value = 0.04381656765025366
$ws.Range("Q2").Value = 6.406819871799779
$ws.Range("R2").Value = 57.66137884619801
$ws.Range("S2").Value = 0.009528950945556487
$ws.Range("T2").Value = 0.009528950945556487
$ws.Range("G3").Value = 27.90240433333333
$ws.Range("H3").Value = 83.707213
$ws.Range("I3").Value = 0.2174736967445081
$ws.Range("J3").Value = 0.2174736967445081
$ws.Range("O3").Value = 0.7547076606638542
$ws.Range("P3").Value = 0.7547076606638543
$ws.Range("Q3").Value = 110.3526884245282
$ws.Range("R3").Value = 993.174195820754
$ws.Range("S3").Value = 0.1641290649259681
$ws.Range("T3").Value = 0.1641290649259682
$ws.Range("G4").Value = 27.90240433333333
$ws.Range("H4").Value = 83.707213
$ws.Range("I4").Value = 0.2174736967445081
$ws.Range("J4").Value = 0.2174736967445081
$ws.Range("M4").Value = 1.055809
$ws.Range("N4").Value = 3.167427
$ws.Range("O4").Value = 0.201475771685892
$ws.Range("P4").Value = 0.2014757716858921
$ws.Range("Q4").Value = 29.45960961677233
$ws.Range("R4").Value = 265.136486550951
$ws.Range("S4").Value = 0.04381568087298343
$ws.Range("T4").Value = 0.04381568087298343
$ws.Range("G5").Value = 64.92210766666668
$ws.Range("I5").Value = 0.5060083921817455
$ws.Range("J5").Value = 0.5060083921817455
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2296153333333334
$ws.Range("N5").Value = 0.6888460000000001
$ws.Range("O5").Value = 0.04381656765025366
$ws.Range("P5").Value = 0.04381656765025366
$ws.Range("Q5").Value = 14.90711139258423
$ws.Range("R5").Value = 134.164002533258
$ws.Range("S5").Value = 0.02217155094762754
$ws.Range("T5").Value = 0.02217155094762754
$ws.Range("G6").Value = 64.92210766666668
$ws.Range("I6").Value = 0.5060083921817455
$ws.Range("J6").Value = 0.5060083921817455
$ws.Range("O6").Value = 0.7547076606638542
$ws.Range("P6").Value = 0.7547076606638543
$ws.Range("Q6").Value = 256.7638628419038
$ws.Range("S6").Value = 0.3818884099397633
$ws.Range("T6").Value = 0.3818884099397633
$ws.Range("G7").Value = 64.92210766666668
$ws.Range("I7").Value = 0.5060083921817455
$ws.Range("J7").Value = 0.5060083921817455
$ws.Range("M7").Value = 1.055809
$ws.Range("N7").Value = 3.167427
$ws.Range("O7").Value = 0.201475771685892
$ws.Range("P7").Value = 0.2014757716858921
$ws.Range("Q7").Value = 68.54534557343568
$ws.Range("R7").Value = 616.9081101609211
$ws.Range("S7").Value = 0.1019484312943547
$ws.Range("T7").Value = 0.1019484312943547
$ws.Range("G8").Value = 19.423329
$ws.Range("H8").Value = 58.269987
$ws.Range("I8").Value = 0.1513870672309258
$ws.Range("J8").Value = 0.1513870672309258
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2296153333333334
$ws.Range("N8").Value = 0.6888460000000001
$ws.Range("O8").Value = 0.04381656765025366
$ws.Range("P8").Value = 0.04381656765025366
$ws.Range("Q8").Value = 4.459894162778
$ws.Range("R8").Value = 40.139047465002
$ws.Range("S8").Value = 0.006633261672697359
$ws.Range("T8").Value = 0.006633261672697361
$ws.Range("G9").Value = 19.423329
$ws.Range("H9").Value = 58.269987
$ws.Range("I9").Value = 0.1513870672309258
$ws.Range("J9").Value = 0.1513870672309258
$ws.Range("O9").Value = 0.7547076606638542
$ws.Range("P9").Value = 0.7547076606638543
$ws.Range("Q9").Value = 76.81834682409399
$ws.Range("R9").Value = 691.365121416846
$ws.Range("S9").Value = 0.1142529793646136
$ws.Range("T9").Value = 0.1142529793646137
$ws.Range("G10").Value = 19.423329
$ws.Range("H10").Value = 58.269987
$ws.Range("I10").Value = 0.1513870672309258
$ws.Range("J10").Value = 0.1513870672309258
$ws.Range("M10").Value = 1.055809
$ws.Range("N10").Value = 3.167427
$ws.Range("O10").Value = 0.201475771685892
$ws.Range("P10").Value = 0.2014757716858921
$ws.Range("Q10").Value = 20.507325568161
$ws.Range("R10").Value = 184.565930113449
$ws.Range("S10").Value = 0.03050082619361479
$ws.Range("T10").Value = 0.0305008261936148
$ws.Range("G11").Value = 16.05459166666667
$ws.Range("H11").Value = 48.163775
$ws.Range("I11").Value = 0.1251308438428206
$ws.Range("J11").Value = 0.1251308438428206
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2296153333333334
$ws.Range("N11").Value = 0.6888460000000001
$ws.Range("O11").Value = 0.04381656765025366
$ws.Range("P11").Value = 0.04381656765025366
$ws.Range("Q11").Value = 3.686380417072223
$ws.Range("R11").Value = 33.17742375365
$ws.Range("S11").Value = 0.005482804084372274
$ws.Range("T11").Value = 0.005482804084372274
$ws.Range("G12").Value = 16.05459166666667
$ws.Range("H12").Value = 48.163775
$ws.Range("I12").Value = 0.1251308438428206
$ws.Range("J12").Value = 0.1251308438428206
$ws.Range("O12").Value = 0.7547076606638542
$ws.Range("P12").Value = 0.7547076606638543
$ws.Range("Q12").Value = 63.49515012432778
$ws.Range("R12").Value = 571.45635111895
$ws.Range("S12").Value = 0.09443720643350914
$ws.Range("T12").Value = 0.09443720643350915
$ws.Range("G13").Value = 16.05459166666667
$ws.Range("H13").Value = 48.163775
$ws.Range("I13").Value = 0.1251308438428206
$ws.Range("J13").Value = 0.1251308438428206
$ws.Range("M13").Value = 1.055809
$ws.Range("N13").Value = 3.167427
$ws.Range("O13").Value = 0.201475771685892
$ws.Range("P13").Value = 0.2014757716858921
$ws.Range("Q13").Value = 16.95058237299167
$ws.Range("R13").Value = 152.555241356925
$ws.Range("S13").Value = 0.02521083332493912
$ws.Range("T13").Value = 0.02521083332493913
